# Scheduled-runner market data refresh: update currentAveragePrice/Leve price/profit columns
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 5662.4136
$ws.Range("I38").Value = 838.7692
$ws.Range("J38").Value = 9581.625
$ws.Range("K38").Value = 2516.3076
$ws.Range("L38").Value = 28744.875
$ws.Range("M38").Value = -2144.3076
$ws.Range("N38").Value = -29488.875

$ws.Range("H43").Value = 3777.889
$ws.Range("I43").Value = 3000
$ws.Range("J43").Value = 4400.2
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 4400.2
$ws.Range("M43").Value = -2931
$ws.Range("N43").Value = -4538.2

$ws.Range("H87").Value = 54500
$ws.Range("I87").Value = 24000
$ws.Range("J87").Value = 58312.5
$ws.Range("K87").Value = 24000
$ws.Range("L87").Value = 58312.5
$ws.Range("M87").Value = -22752
$ws.Range("N87").Value = -60808.5

$ws.Range("H90").Value = 54500
$ws.Range("I90").Value = 24000
$ws.Range("J90").Value = 58312.5
$ws.Range("K90").Value = 72000
$ws.Range("L90").Value = 174937.5
$ws.Range("M90").Value = -65760
$ws.Range("N90").Value = -187417.5

$ws.Range("H121").Value = 3000
$ws.Range("I121").Value = 3000
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 9000
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = -7253
$ws.Range("N121").ClearContents()

$ws.Range("H132").Value = 2391
$ws.Range("I132").Value = 2388.6
$ws.Range("J132").Value = 2415
$ws.Range("K132").Value = 7165.799999999999
$ws.Range("L132").Value = 7245
$ws.Range("M132").Value = -4635.799999999999

$ws.Range("H138").Value = 2445.0195
$ws.Range("I138").Value = 1885.0541
$ws.Range("J138").Value = 3924.9285
$ws.Range("K138").Value = 5655.1623
$ws.Range("L138").Value = 11774.7855
$ws.Range("M138").Value = -515.1623
$ws.Range("N138").Value = -22054.7855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 89354
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 89354
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 89354
$ws.Range("N24").Value = -90102

$ws.Range("H32").Value = 2256.4058
$ws.Range("I32").Value = 1517.9539
$ws.Range("J32").Value = 14256.25
$ws.Range("K32").Value = 1517.9539
$ws.Range("L32").Value = 14256.25
$ws.Range("M32").Value = -1230.9539

$ws.Range("H61").Value = 3472.625
$ws.Range("I61").Value = 1963.8334
$ws.Range("J61").Value = 7999
$ws.Range("K61").Value = 1963.8334
$ws.Range("L61").Value = 7999
$ws.Range("M61").Value = -1751.8334
$ws.Range("N61").Value = -8423

$ws.Range("H74").Value = 4816
$ws.Range("I74").Value = 4613.7856
$ws.Range("J74").Value = 5099.1
$ws.Range("K74").Value = 4613.7856
$ws.Range("L74").Value = 5099.1
$ws.Range("M74").Value = -3739.7856

$ws.Range("H77").Value = 4816
$ws.Range("I77").Value = 4613.7856
$ws.Range("J77").Value = 5099.1
$ws.Range("K77").Value = 23068.928
$ws.Range("L77").Value = 25495.5
$ws.Range("M77").Value = -18700.928

$ws.Range("H88").Value = 4543.4
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 4543.4
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 4543.4
$ws.Range("N88").Value = -5355.4

$ws.Range("H91").Value = 4543.4
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 4543.4
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 4543.4
$ws.Range("N91").Value = -7351.4

$ws.Range("H100").Value = 89354
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 89354
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 89354
$ws.Range("N100").Value = -91518

$ws.Range("H132").Value = 10772.228
$ws.Range("I132").Value = 5383.8203
$ws.Range("J132").Value = 52801.8
$ws.Range("K132").Value = 16151.4609
$ws.Range("L132").Value = 158405.4
$ws.Range("M132").Value = -13621.4609

$ws.Range("H136").Value = 3472.625
$ws.Range("I136").Value = 1963.8334
$ws.Range("J136").Value = 7999
$ws.Range("K136").Value = 5891.5002
$ws.Range("L136").Value = 23997
$ws.Range("M136").Value = -3341.5002
$ws.Range("N136").Value = -29097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 929.48334
$ws.Range("I94").Value = 644.5172
$ws.Range("J94").Value = 1196.0646
$ws.Range("K94").Value = 644.5172
$ws.Range("L94").Value = 1196.0646
$ws.Range("M94").Value = -193.5172

$ws.Range("H134").Value = 3839.3794
$ws.Range("I134").Value = 3394.2104
$ws.Range("J134").Value = 4685.2
$ws.Range("K134").Value = 10182.6312
$ws.Range("L134").Value = 14055.6
$ws.Range("M134").Value = -7647.6312

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 807.5172
$ws.Range("I22").Value = 324.22223
$ws.Range("J22").Value = 1598.3636
$ws.Range("K22").Value = 324.22223
$ws.Range("L22").Value = 1598.3636
$ws.Range("M22").Value = 25.77776999999998
$ws.Range("N22").Value = -2298.3636

$ws.Range("H58").Value = 4249.7856
$ws.Range("I58").Value = 4632.0835
$ws.Range("J58").Value = 1956
$ws.Range("K58").Value = 4632.0835
$ws.Range("L58").Value = 1956
$ws.Range("M58").Value = -4429.0835

$ws.Range("H62").Value = 10026.611
$ws.Range("I62").Value = 7618.7
$ws.Range("J62").Value = 13036.5
$ws.Range("K62").Value = 7618.7
$ws.Range("L62").Value = 13036.5
$ws.Range("M62").Value = -6994.7
$ws.Range("N62").Value = -14284.5

$ws.Range("H65").Value = 10026.611
$ws.Range("I65").Value = 7618.7
$ws.Range("J65").Value = 13036.5
$ws.Range("K65").Value = 38093.5
$ws.Range("L65").Value = 65182.5
$ws.Range("M65").Value = -34973.5
$ws.Range("N65").Value = -71422.5

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H134").Value = 2194.9473
$ws.Range("I134").Value = 2009.3846
$ws.Range("J134").Value = 2597
$ws.Range("K134").Value = 6028.1538
$ws.Range("L134").Value = 7791
$ws.Range("M134").Value = -3493.1538
$ws.Range("N134").Value = -12861

$ws.Range("H136").Value = 4249.7856
$ws.Range("I136").Value = 4632.0835
$ws.Range("J136").Value = 1956
$ws.Range("K136").Value = 13896.2505
$ws.Range("L136").Value = 5868
$ws.Range("M136").Value = -11346.2505

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 58830720
$ws.Range("I4").Value = 66666816
$ws.Range("J4").Value = 60005
$ws.Range("K4").Value = 200000448
$ws.Range("L4").Value = 180015
$ws.Range("M4").Value = -200000336
$ws.Range("N4").Value = -180239

$ws.Range("H37").Value = 105440.336
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 105440.336
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 316321.008
$ws.Range("N37").Value = -316545.008

$ws.Range("H39").Value = 5628.5884
$ws.Range("I39").Value = 4000
$ws.Range("J39").Value = 6307.1665
$ws.Range("K39").Value = 12000
$ws.Range("L39").Value = 18921.4995
$ws.Range("M39").Value = -11706
$ws.Range("N39").Value = -19509.4995

$ws.Range("H60").Value = 2015

$ws.Range("H132").Value = 1750.6

$ws.Range("H141").Value = 4723.222
$ws.Range("I141").Value = 1877.5
$ws.Range("J141").Value = 6999.8
$ws.Range("K141").Value = 5632.5
$ws.Range("L141").Value = 20999.4
$ws.Range("M141").Value = -452.5
$ws.Range("N141").Value = -31359.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()

$ws.Range("H94").Value = 59999
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 59999
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 59999
$ws.Range("N94").Value = -61351

$ws.Range("H107").Value = 1292.45
$ws.Range("I107").Value = 1927.1
$ws.Range("J107").Value = 657.8
$ws.Range("K107").Value = 1927.1
$ws.Range("L107").Value = 657.8
$ws.Range("M107").Value = -7.099999999999909
$ws.Range("N107").Value = -4497.8

$ws.Range("H132").Value = 3320.4119
$ws.Range("I132").Value = 2917.9614
$ws.Range("J132").Value = 4628.375
$ws.Range("K132").Value = 8753.8842
$ws.Range("L132").Value = 13885.125
$ws.Range("M132").Value = -6223.8842
$ws.Range("N132").Value = -18945.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2018.5652
$ws.Range("I22").Value = 657.7692
$ws.Range("J22").Value = 3787.6
$ws.Range("K22").Value = 657.7692
$ws.Range("L22").Value = 3787.6
$ws.Range("M22").Value = -362.7692

$ws.Range("H27").Value = 2018.5652
$ws.Range("I27").Value = 657.7692
$ws.Range("J27").Value = 3787.6
$ws.Range("K27").Value = 657.7692
$ws.Range("L27").Value = 3787.6
$ws.Range("M27").Value = -550.7692

$ws.Range("H132").Value = 5421.5386
$ws.Range("I132").Value = 4448.3
$ws.Range("J132").Value = 8665.666999999999
$ws.Range("K132").Value = 13344.9
$ws.Range("L132").Value = 25997.001
$ws.Range("M132").Value = -10814.9

$ws.Range("H136").Value = 5324.769
$ws.Range("I136").Value = 4843.5713
$ws.Range("J136").Value = 7345.8
$ws.Range("K136").Value = 14530.7139
$ws.Range("L136").Value = 22037.4
$ws.Range("M136").Value = -11980.7139
$ws.Range("N136").Value = -27137.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1965.1666
$ws.Range("I126").Value = 2058.9
$ws.Range("J126").Value = 1496.5
$ws.Range("K126").Value = 6176.700000000001
$ws.Range("L126").Value = 4489.5
$ws.Range("M126").Value = -3706.700000000001
$ws.Range("N126").Value = -9429.5

$ws.Range("H132").Value = 2226.3157
$ws.Range("I132").Value = 1437.7273
$ws.Range("J132").Value = 3310.625
$ws.Range("K132").Value = 4313.1819
$ws.Range("L132").Value = 9931.875
$ws.Range("M132").Value = -1783.1819
